# Generate Report for Handoff
# Update the "Latest Handoff Date/Datetime" timestamps on all three sheets
# for every file row that is being (re-)handed off in this report run
# (i.e. every row except the ones already "Handed back: in sync with en-US"
# and the one still "In Translation").

$wb = $excel.ActiveWorkbook

# --- Sheet "Overview": column D = "Latest Handoff Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D4").Value  = "2016-03-24 11:11:42"
$wsOverview.Range("D6").Value  = "2016-03-24 11:11:42"
$wsOverview.Range("D7").Value  = "2016-03-24 11:11:42"
$wsOverview.Range("D8").Value  = "2016-03-24 11:11:42"
$wsOverview.Range("D9").Value  = "2016-03-24 11:11:42"
$wsOverview.Range("D10").Value = "2016-03-24 11:11:42"

# --- Sheet "zh-cn": column E = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value  = "2016-03-24 11:11:31"
$wsZhCn.Range("E6").Value  = "2016-03-24 11:11:31"
$wsZhCn.Range("E7").Value  = "2016-03-24 11:11:31"
$wsZhCn.Range("E8").Value  = "2016-03-24 11:11:31"
$wsZhCn.Range("E9").Value  = "2016-03-24 11:11:31"
$wsZhCn.Range("E10").Value = "2016-03-24 11:11:31"

# --- Sheet "de-de": column E = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value  = "2016-03-24 11:11:42"
$wsDeDe.Range("E6").Value  = "2016-03-24 11:11:42"
$wsDeDe.Range("E7").Value  = "2016-03-24 11:11:42"
$wsDeDe.Range("E8").Value  = "2016-03-24 11:11:42"
$wsDeDe.Range("E9").Value  = "2016-03-24 11:11:42"
$wsDeDe.Range("E10").Value = "2016-03-24 11:11:42"
